# planillaZulliger.xlsx - "corrigiendo valores del excel"
#
# 1. Shorten two header labels.
# 2. Turn the "Lam" (A) column into real numbers instead of text, and
#    fix/replace a bunch of "?" placeholder values across rows 2-4.
# 3. Fold the lone row that lived on the extra "Sheet1" tab back into
#    "Hoja de datos" as new rows 5-7 (one row per "Lam" 1/2/3, mirroring
#    the pattern already used by rows 2-4), then correct their values.
# 4. Delete the now-redundant "Sheet1" tab.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja de datos")
$extra = $wb.Worksheets.Item("Sheet1")

# --- Header row tweaks ------------------------------------------------------
$ws.Range("F1").Value = "Det"
$ws.Range("I1").Value = "Cont"

# --- Row 2 -------------------------------------------------------------------
$ws.Range("A2").Value = 1
$ws.Range("F2").Value = "FM"
$ws.Range("H2").Value = ""
$ws.Range("I2").Value = "A"

# --- Row 3 -------------------------------------------------------------------
# A3 currently holds the text "2" - clone it (as text) into H3 before turning
# A3 itself into a real number, so H3 keeps a *text* "2".
$ws.Range("A3").Copy()
$ws.Range("H3").PasteSpecial()
$ws.Range("A3").Value = 2
$ws.Range("F3").Value = "C',M"
$ws.Range("I3").Value = "H,H"

# --- Row 4 -------------------------------------------------------------------
$ws.Range("A4").Value = 3
$ws.Range("F4").Value = "m"
$ws.Range("H4").Value = ""
$ws.Range("I4").Value = "Fi"
$ws.Range("J4").Value = ""

# --- Rows 5-7: seed from the old Sheet1 row (keeps its text cell types, e.g.
#     B/C/D/E/G/K/L stay text "1"/"?" instead of turning into numbers), then
#     correct the per-row values -------------------------------------------
$extra.Range("A5:L5").Copy()
$ws.Range("A5").PasteSpecial()
$extra.Range("A5:L5").Copy()
$ws.Range("A6").PasteSpecial()
$extra.Range("A5:L5").Copy()
$ws.Range("A7").PasteSpecial()

# Row 5
$ws.Range("A5").Value = 1
$ws.Range("F5").Value = ""
$ws.Range("H5").Value = ""
$ws.Range("I5").Value = "Ad,Ad"
$ws.Range("J5").Value = ""

# Row 6
$ws.Range("A6").Value = 2
$ws.Range("F6").Value = ""
$ws.Range("H6").Value = ""
$ws.Range("I6").Value = ""
$ws.Range("J6").Value = ""

# Row 7
$ws.Range("A7").Value = 3
$ws.Range("F7").Value = "C"
$ws.Range("H7").Value = ""
$ws.Range("I7").Value = "Hx"
$ws.Range("J7").Value = ""

# --- Remove the now-redundant extra worksheet --------------------------------
$excel.DisplayAlerts = $false
[void]$extra.Delete()
$excel.DisplayAlerts = $true
